$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the run-content of a paragraph (everything except the
# trailing paragraph mark) with an explicit sequence of <w:r>...</w:r> runs
# supplied as raw OOXML. Using InsertXML (rather than Find/Replace) keeps
# each <w:r> as a *separate* run instead of letting adjacent, identically
# formatted runs be silently recombined into one.
#
# NOTE: this runtime's PowerShell only reliably marshals COM objects through
# *positional* parameters - named parameters (-Foo bar) silently drop COM
# values - so every helper below is called positionally.
# ---------------------------------------------------------------------------
function Set-ParagraphRuns($Paragraph, $RunsXml) {
    $pRange = $Paragraph.Range
    $target = $d.Range($pRange.Start, $pRange.End - 1)

    $pkg = '<?xml version="1.0" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $RunsXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($pkg)
}

function New-Run($RprXml, $Text, $Preserve) {
    $sp = ""
    if ($Preserve) { $sp = ' xml:space="preserve"' }
    return "<w:r>$RprXml<w:t$sp>$Text</w:t></w:r>"
}

function New-RunBreak($RprXml) {
    return "<w:r>$RprXml<w:br/></w:r>"
}

function New-RunBreakText($RprXml, $Text, $Preserve) {
    $sp = ""
    if ($Preserve) { $sp = ' xml:space="preserve"' }
    return "<w:r>$RprXml<w:br/><w:t$sp>$Text</w:t></w:r>"
}

function New-RunPageBreakText($RprXml, $Text, $Preserve) {
    $sp = ""
    if ($Preserve) { $sp = ' xml:space="preserve"' }
    return "<w:r>$RprXml<w:lastRenderedPageBreak/><w:t$sp>$Text</w:t></w:r>"
}

# rPr blocks reused throughout the document
$rpr44 = '<w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="44"/></w:rPr>'
$rpr36 = '<w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr>'
$rpr32 = '<w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr>'
$rpr24 = '<w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr>'
$rpr28 = '<w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="28"/></w:rPr>'
$rprN  = '<w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr>'

# ---------------------------------------------------------------------------
# Paragraph 1: Title
# ---------------------------------------------------------------------------
$runs1 = New-Run $rpr44 "Exploring the Nexus of Biology, Medicine, and Society" $false
Set-ParagraphRuns $d.Paragraphs.Item(1) $runs1

# ---------------------------------------------------------------------------
# Paragraph 2: Author name -> "Dr. Evelyn Graham, Ph.D." split over 6 runs
# ---------------------------------------------------------------------------
$runs2 = New-Run $rpr36 "Dr" $false
$runs2 += New-Run $rpr36 "." $false
$runs2 += New-Run $rpr36 " Evelyn Graham, Ph" $true
$runs2 += New-Run $rpr36 "." $false
$runs2 += New-Run $rpr36 "D" $false
$runs2 += New-Run $rpr36 "." $false
Set-ParagraphRuns $d.Paragraphs.Item(2) $runs2

# ---------------------------------------------------------------------------
# Paragraph 3: email -> evelyn.g@riley.edu
# ---------------------------------------------------------------------------
$runs3 = New-Run $rpr32 "evelyn" $false
$runs3 += New-Run $rpr32 "." $false
$runs3 += New-Run $rpr32 "g@riley" $false
$runs3 += New-Run $rpr32 "." $false
$runs3 += New-Run $rpr32 "edu" $false
Set-ParagraphRuns $d.Paragraphs.Item(3) $runs3

# Paragraph 4 stays an empty paragraph - no change needed.

# ---------------------------------------------------------------------------
# Paragraph 5: main body text
# ---------------------------------------------------------------------------
$runs5  = New-Run $rpr24 "Journey into the intricate realm where biology and medicine intertwine with the tapestry of human existence" $false
$runs5 += New-Run $rpr24 "." $false
$runs5 += New-Run $rpr24 " Biology unveils the marvelous intricacies of life, delving into the cellular mysteries that orchestrate the symphony of life" $true
$runs5 += New-Run $rpr24 "." $false
$runs5 += New-Run $rpr24 " Medicine, its counterpart, embarks on a noble quest to alleviate suffering and restore health, employing scientific knowledge to heal, prevent, and cure ailments" $true
$runs5 += New-Run $rpr24 "." $false
$runs5 += New-Run $rpr24 " Amidst this dynamic interplay, society plays a pivotal role, shaping and being shaped by advancements in biological and medical understanding" $true
$runs5 += New-Run $rpr24 "." $false
$runs5 += New-Run $rpr24 " As we embark on this scholastic odyssey, we shall explore the profound impact of biology and medicine on societal progress, delving into case studies that illuminate the symbiotic relationship between these disciplines and the communities they serve" $true
$runs5 += New-Run $rpr24 "." $false
$runs5 += New-RunBreak $rpr24
$runs5 += New-RunBreakText $rpr24 "Exploring the profound impact of biology and medicine on the trajectory of human civilization unveils a saga of remarkable achievements that have transformed our understanding of life and health" $false
$runs5 += New-Run $rpr24 "." $false
$runs5 += New-Run $rpr24 " From the advent of antibiotics to the advent of vaccines, scientific breakthroughs in these fields have led to a dramatic reduction in infectious diseases, extending human life expectancy and improving overall well-being" $true
$runs5 += New-Run $rpr24 "." $false
$runs5 += New-Run $rpr24 " The Human Genome Project stands as a testament to scientific ingenuity, paving the way for personalized medicine, targeted therapies, and a deeper comprehension of genetic disorders" $true
$runs5 += New-Run $rpr24 "." $false
$runs5 += New-Run $rpr24 " Moreover, the field of genomics continues to revolutionize our insights into evolution, biodiversity, and conservation efforts" $true
$runs5 += New-Run $rpr24 "." $false
$runs5 += New-RunBreak $rpr24
$runs5 += New-RunBreakText $rpr24 "The connection between biology, medicine, and society is a dynamic interplay, marked by both challenges and opportunities" $false
$runs5 += New-Run $rpr24 "." $false
$runs5 += New-Run $rpr24 " The rise of antimicrobial resistance poses a significant public health threat, demanding the development of novel antimicrobial agents and prudent antibiotic stewardship" $true
$runs5 += New-Run $rpr24 "." $false
$runs5 += New-Run $rpr24 " Additionally, the ethical, legal, and social implications of genetic testing and gene editing necessitate careful consideration as we navigate the complexities of these evolving technologies" $true
$runs5 += New-Run $rpr24 "." $false
$runs5 += New-Run $rpr24 " Recognizing the profound influence of biology and medicine on societal progress, it " $true
$runs5 += New-RunPageBreakText $rpr24 "becomes imperative to foster interdisciplinary collaborations, promote scientific literacy, and ensure equitable access to healthcare services" $false
$runs5 += New-Run $rpr24 "." $false
Set-ParagraphRuns $d.Paragraphs.Item(5) $runs5

# Paragraph 6 "Summary" heading is unchanged.

# ---------------------------------------------------------------------------
# Paragraph 7: summary body
# ---------------------------------------------------------------------------
$runs7  = New-Run $rprN "Through this comprehensive journey, we have delved into the profound intertwined relationship between biology, medicine, and society" $false
$runs7 += New-Run $rprN "." $false
$runs7 += New-Run $rprN " From the remarkable triumphs of disease control and genetic discoveries to the ongoing challenges of antimicrobial resistance and ethical considerations, this exploration has illuminated the transformative impact of these disciplines on human well-being" $true
$runs7 += New-Run $rprN "." $false
$runs7 += New-Run $rprN " As we continue to advance our understanding of life and health, it becomes increasingly evident that biology, medicine, and society are inextricably linked, shaping and being shaped by each other in a continuous cycle of progress and challenge" $true
$runs7 += New-Run $rprN "." $false
Set-ParagraphRuns $d.Paragraphs.Item(7) $runs7

# ---------------------------------------------------------------------------
# Add a new empty paragraph at the very end of the document body (before the
# final section properties), mirroring the appended <w:p/> in the diff.
# ---------------------------------------------------------------------------
$endOfDoc = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endOfDoc.InsertParagraphAfter()
